$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.901.82"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.633.55"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.56"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0882"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.864.59"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.649.83"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.43"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "27.901.88"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.21"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.08"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.51"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "1.393.25"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  +10.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.560"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.850"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.83"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.91"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.44"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").Value = "1.773.27"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.81"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.67"
$ws.Range("E51").Value = "  +1.40%  "
